# Apply sprint-by-sprint data to the "User Stories" Planning Poker workbook.
# Adds a "Sprint" column (D), a per-sprint "Sprint Story Point Total" column
# (E) with subtotal formulas, and Total/Average summary rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Stories")

$xlCenter = -4108

# --- Sprint number column (D), grouped by the same colour blocks already
#     used in column A: Sprint1=green rows 2-13, Sprint2=red rows 14-25,
#     Sprint3=blue rows 26-37, Sprint4=purple rows 38-49. ---
$ws.Range("D2:D13").Value = 1
$ws.Range("D2:D13").Font.Color = $ws.Range("A2").Font.Color

$ws.Range("D14:D25").Value = 2
$ws.Range("D14:D25").Font.Color = $ws.Range("A14").Font.Color

$ws.Range("D26:D37").Value = 3
$ws.Range("D26:D37").Font.Color = $ws.Range("A26").Font.Color

# --- Per-sprint story-point subtotal formulas in column E (row = last row
#     of each sprint block minus one, matching the source layout). ---
$ws.Range("E31").Formula = "=SUM(B26:B37)"
$ws.Range("E31").Font.Color = $ws.Range("A26").Font.Color
$ws.Range("E31").HorizontalAlignment = $xlCenter

$ws.Range("E7").Formula = "=SUM(B2:B13)"
$ws.Range("E7").Font.Color = $ws.Range("A2").Font.Color
$ws.Range("E7").HorizontalAlignment = $xlCenter

$ws.Range("E19").Formula = "=SUM(B14:B25)"
$ws.Range("E19").Font.Color = $ws.Range("A14").Font.Color
$ws.Range("E19").HorizontalAlignment = $xlCenter

# --- Sprint 4 (rows 38-49): recolour the story text purple and add the
#     sprint number in column D, also purple. ---
$purple = 0xA03070
$ws.Range("A38:A49").Font.Color = $purple
$ws.Range("D38:D49").Value = 4
$ws.Range("D38:D49").Font.Color = $purple

$ws.Range("E43").Formula = "=SUM(B38:B49)"
$ws.Range("E43").Font.Color = $purple
$ws.Range("E43").HorizontalAlignment = $xlCenter

# --- Header row (row 1): bold, plus the two new column headers. ---
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("D1").Value = "Sprint "
$ws.Range("D1").Font.Bold = $true
$ws.Range("E1").Value = "Sprint Story Point Total"
$ws.Range("E1").Font.Bold = $true

# --- Totals / averages summary rows. ---
$ws.Range("A50").Value = "Total Story Points:"
$ws.Range("A50").Font.Bold = $true
$ws.Range("B50").Formula = "=SUM(B2:B49)"
$ws.Range("C50").Formula = "=SUM(C2:C49)"

$ws.Range("A51").Value = "Average Story Points:"
$ws.Range("A51").Font.Bold = $true
$ws.Range("B51").Formula = "=AVERAGE(B2:B49)"
$ws.Range("C51").Formula = "=AVERAGE(C2:C49)"

# --- Blank styled row at the bottom (matches the source colour block). ---
$ws.Range("A53").Font.Color = $ws.Range("A2").Font.Color

# --- New column E width. ---
$ws.Columns.Item(5).ColumnWidth = 18.29

# --- Scroll the sheet view back to the top. ---
$ws.Application.ActiveWindow.ScrollRow = 1

$wb.Save()
